$d = $word.ActiveDocument

# Document order of the 7 rectangle/Line shapes and their target stroke
# widths (in points == a:ln/@w in EMU, 1pt = 12700 EMU):
#   1) lgDashDotDot (Téglalap 7)  50800 -> 76200 (6pt)
#   2) lgDashDot    (Téglalap 6)  50800 -> 38100 (3pt)
#   3) dashDot      (Téglalap 4)  50800 -> 76200 (6pt)
#   4) dash         (Téglalap 3)  50800 -> 57150 (4.5pt)
#   5) sysDash      (Téglalap 2)  50800 -> 28575 (2.25pt)
#   6) lgDash       (Téglalap 5)  50800 -> 19050 (1.5pt)
#   7) sysDot       (Rectangle 1) 50800 -> 12700 (1pt)
$d.Shapes.Item(1).Line.Weight = 6
$d.Shapes.Item(2).Line.Weight = 3
$d.Shapes.Item(3).Line.Weight = 6
$d.Shapes.Item(4).Line.Weight = 4.5
$d.Shapes.Item(5).Line.Weight = 2.25
$d.Shapes.Item(6).Line.Weight = 1.5
$d.Shapes.Item(7).Line.Weight = 1
